$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (shifts old J..K to K..L), carrying
# over cell values/styles and the custom column width.
[void]$ws.Range("J1").EntireColumn.Insert()

# Grow the table ("Tableau1") to include the freshly inserted column.
$tbl = $ws.ListObjects.Item("Tableau1")
[void]$tbl.Resize($ws.Range("A1:K7"))

# Set the header text for the new column and the (now shifted) neighbour so
# the table's ListColumn names re-sync from the header row.
$ws.Range("J1").Value = "Close to"
$ws.Range("K1").Value = "is weekend ok?"

# Data entered in the new "Close to" column.
$ws.Range("J2").Value = "TP424"

# The two existing data-validation rules reference column K, which has now
# shifted to column L - repoint them (their sqref already auto-shifted).
[void]$ws.Range("E2:E7").Validation.Modify(3, 1, 1, "=`$L`$4:`$L`$5")
[void]$ws.Range("K2:K7").Validation.Modify(3, 1, 1, "=`$L`$2:`$L`$3")

[void]$ws.Range("J3").Select()
